$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 59917706715738.73
$ws.Range("C2").Value = 59917706715738.74
$ws.Range("D2").Value = 59917706715738.74

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 41582531770713.63
$ws.Range("C3").Value = 41610723317676.81
$ws.Range("D3").Value = 46177753925714.51

# Row 4 - name change GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 51001528979309.88
$ws.Range("C4").Value = 51129895334717.09
$ws.Range("D4").Value = 49340241390407.26

# Row 5 - name change AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 79099859590571.61
$ws.Range("C5").Value = 52788620692140.96
$ws.Range("D5").Value = 69048803814534.45
